$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 2638.4146
$ws.Range("I15").Value = 2638.4146
$ws.Range("K15").Value = 7915.2438
$ws.Range("M15").Value = -7746.2438
# Row 19
$ws.Range("H19").Value = 415.7143
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 415.7143
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 415.7143
$ws.Range("N19").Value = -765.7143
$ws.Range("M19").ClearContents()
# Row 69
$ws.Range("H69").Value = 27071.143
$ws.Range("J69").Value = 15199.8
$ws.Range("L69").Value = 45599.39999999999
$ws.Range("N69").Value = -47347.39999999999
# Row 72
$ws.Range("H72").Value = 27071.143
$ws.Range("J72").Value = 15199.8
$ws.Range("L72").Value = 136798.2
$ws.Range("N72").Value = -145534.2
# Row 80
$ws.Range("H80").Value = 587.8
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 587.8
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 1763.4
$ws.Range("N80").Value = -3759.4
$ws.Range("M80").ClearContents()
# Row 83
$ws.Range("H83").Value = 587.8
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 587.8
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 5290.2
$ws.Range("N83").Value = -15274.2
$ws.Range("M83").ClearContents()
# Row 112
$ws.Range("H112").Value = 3933.611
$ws.Range("I112").Value = 1437.5
$ws.Range("J112").Value = 4646.7856
$ws.Range("K112").Value = 4312.5
$ws.Range("L112").Value = 13940.3568
$ws.Range("M112").Value = -3204.5
$ws.Range("N112").Value = -16156.3568
# Row 116
$ws.Range("H116").Value = 6241
$ws.Range("J116").Value = 5306.857
$ws.Range("L116").Value = 5306.857
$ws.Range("N116").Value = -12190.857
# Row 121
$ws.Range("H121").Value = 3933
$ws.Range("J121").Value = 3933
$ws.Range("L121").Value = 11799
$ws.Range("N121").Value = -15293
# Row 132
$ws.Range("H132").Value = 10038.182
$ws.Range("I132").Value = 7738.1724
$ws.Range("J132").Value = 26713.25
$ws.Range("K132").Value = 23214.5172
$ws.Range("L132").Value = 80139.75
$ws.Range("M132").Value = -20684.5172
$ws.Range("N132").Value = -85199.75
# Row 137
$ws.Range("H137").Value = 8122.4165
$ws.Range("I137").Value = 2968.353
$ws.Range("J137").Value = 12733.947
$ws.Range("K137").Value = 8905.059000000001
$ws.Range("L137").Value = 38201.841
$ws.Range("M137").Value = -6355.059000000001
$ws.Range("N137").Value = -43301.841
# Row 141
$ws.Range("H141").Value = 2886.9395
$ws.Range("I141").Value = 2469.9678
$ws.Range("K141").Value = 7409.903399999999
$ws.Range("M141").Value = -2229.903399999999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 6675.2856
$ws.Range("I2").Value = 5710.84
$ws.Range("J2").Value = 9086.4
$ws.Range("K2").Value = 5710.84
$ws.Range("L2").Value = 9086.4
$ws.Range("M2").Value = -5597.84
$ws.Range("N2").Value = -9312.4
# Row 24
$ws.Range("H24").Value = 27677.5
$ws.Range("J24").Value = 27677.5
$ws.Range("L24").Value = 27677.5
$ws.Range("N24").Value = -28425.5
# Row 43
$ws.Range("H43").Value = 35000
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 35000
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 35000
$ws.Range("N43").Value = -35626
$ws.Range("M43").ClearContents()
# Row 61
$ws.Range("H61").Value = 85431.42999999999
$ws.Range("I61").Value = 2346.25
$ws.Range("J61").Value = 196211.67
$ws.Range("K61").Value = 2346.25
$ws.Range("L61").Value = 196211.67
$ws.Range("M61").Value = -2134.25
$ws.Range("N61").Value = -196635.67
# Row 100
$ws.Range("H100").Value = 27677.5
$ws.Range("J100").Value = 27677.5
$ws.Range("L100").Value = 27677.5
$ws.Range("N100").Value = -29841.5
# Row 116
$ws.Range("H116").Value = 6675.2856
$ws.Range("I116").Value = 5710.84
$ws.Range("J116").Value = 9086.4
$ws.Range("K116").Value = 5710.84
$ws.Range("L116").Value = 9086.4
$ws.Range("M116").Value = -3416.84
$ws.Range("N116").Value = -13674.4
# Row 132
$ws.Range("H132").Value = 12270.096
$ws.Range("I132").Value = 1504.8518
$ws.Range("K132").Value = 4514.555399999999
$ws.Range("M132").Value = -1984.555399999999
# Row 136
$ws.Range("H136").Value = 85431.42999999999
$ws.Range("I136").Value = 2346.25
$ws.Range("J136").Value = 196211.67
$ws.Range("K136").Value = 7038.75
$ws.Range("L136").Value = 588635.01
$ws.Range("M136").Value = -4488.75
$ws.Range("N136").Value = -593735.01

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 6675.2856
$ws.Range("I3").Value = 5710.84
$ws.Range("J3").Value = 9086.4
$ws.Range("K3").Value = 5710.84
$ws.Range("L3").Value = 9086.4
$ws.Range("M3").Value = -5596.84
$ws.Range("N3").Value = -9314.4
# Row 63
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
# Row 66
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
# Row 99
$ws.Range("H99").Value = 19337.793
$ws.Range("I99").Value = 22058.95
$ws.Range("J99").Value = 13290.777
$ws.Range("K99").Value = 22058.95
$ws.Range("L99").Value = 13290.777
$ws.Range("M99").Value = -20560.95
$ws.Range("N99").Value = -16286.777

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 979.56525
$ws.Range("I22").Value = 616.86664
$ws.Range("K22").Value = 616.86664
$ws.Range("M22").Value = -266.86664
# Row 76
$ws.Range("H76").Value = 8437
$ws.Range("I76").Value = 8437
$ws.Range("K76").Value = 8437
$ws.Range("M76").Value = -8122
# Row 79
$ws.Range("H79").Value = 8437
$ws.Range("I79").Value = 8437
$ws.Range("K79").Value = 8437
$ws.Range("M79").Value = -7345
# Row 125
$ws.Range("H125").Value = 76331.664
$ws.Range("J125").Value = 76331.664
$ws.Range("L125").Value = 76331.664
$ws.Range("N125").Value = -81251.664
# Row 134
$ws.Range("H134").Value = 4568.4565
$ws.Range("I134").Value = 1551.8206
$ws.Range("J134").Value = 21375.428
$ws.Range("K134").Value = 4655.4618
$ws.Range("L134").Value = 64126.284
$ws.Range("M134").Value = -2120.4618
$ws.Range("N134").Value = -69196.284

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 334.16666
$ws.Range("I2").Value = 406.77777
$ws.Range("J2").Value = 116.333336
$ws.Range("K2").Value = 2440.66662
$ws.Range("L2").Value = 698.000016
$ws.Range("M2").Value = -2327.66662
$ws.Range("N2").Value = -924.000016
# Row 97
$ws.Range("H97").Value = 320.84616
$ws.Range("I97").Value = 270.0909
$ws.Range("K97").Value = 810.2727
$ws.Range("M97").Value = -314.2727
# Row 107
$ws.Range("H107").Value = 8164.3706
$ws.Range("I107").Value = 17222.666
$ws.Range("K107").Value = 51667.99800000001
$ws.Range("M107").Value = -49747.99800000001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 14822.148
$ws.Range("I132").Value = 12905.333
$ws.Range("K132").Value = 38715.999
$ws.Range("M132").Value = -36185.999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 11
$ws.Range("H11").Value = 15426.625
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 15426.625
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 15426.625
$ws.Range("N11").Value = -15706.625
$ws.Range("M11").ClearContents()
# Row 55
$ws.Range("H55").Value = 1078.4117
$ws.Range("I55").Value = 636.3043
$ws.Range("K55").Value = 636.3043
$ws.Range("M55").Value = -463.3043
# Row 122
$ws.Range("H122").Value = 6363.952
$ws.Range("I122").Value = 5429.2964
$ws.Range("K122").Value = 16287.8892
$ws.Range("M122").Value = -13837.8892
# Row 132
$ws.Range("H132").Value = 26960.8
$ws.Range("I132").Value = 16351.333
$ws.Range("J132").Value = 31507.715
$ws.Range("K132").Value = 49053.999
$ws.Range("L132").Value = 94523.145
$ws.Range("M132").Value = -46523.999
$ws.Range("N132").Value = -99583.145

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 28412.875
$ws.Range("I62").Value = 26852.572
$ws.Range("J62").Value = 29055.354
$ws.Range("K62").Value = 26852.572
$ws.Range("L62").Value = 29055.354
$ws.Range("M62").Value = -26228.572
$ws.Range("N62").Value = -30303.354
# Row 65
$ws.Range("H65").Value = 28412.875
$ws.Range("I65").Value = 26852.572
$ws.Range("J65").Value = 29055.354
$ws.Range("K65").Value = 134262.86
$ws.Range("L65").Value = 145276.77
$ws.Range("M65").Value = -131142.86
$ws.Range("N65").Value = -151516.77
# Row 97
$ws.Range("H97").Value = 38717.25
$ws.Range("J97").Value = 38717.25
$ws.Range("L97").Value = 38717.25
$ws.Range("N97").Value = -40699.25
# Row 100
$ws.Range("H100").Value = 659.64703
$ws.Range("J100").Value = 879
$ws.Range("L100").Value = 1758
$ws.Range("N100").Value = -2840
# Row 101
$ws.Range("H101").Value = 26358.334
$ws.Range("J101").Value = 23630
$ws.Range("L101").Value = 23630
$ws.Range("N101").Value = -30120
# Row 132
$ws.Range("H132").Value = 7452.8335
$ws.Range("I132").Value = 2529.6
$ws.Range("K132").Value = 7588.799999999999
$ws.Range("M132").Value = -5058.799999999999
